$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.757.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.16%  "
$ws.Range("D3").Value = "'2.443.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'575.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").Value = "'146.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.76%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.539"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").Value = "'2.442.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("E10").Value = "  +4.98%  "
$ws.Range("D11").Value = "'0.159"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "'5.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "'0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "'27.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.74%  "
$ws.Range("D15").Value = "'0.0000179"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.83%  "
$ws.Range("D16").Value = "'2.889.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").Value = "'62.532.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.35%  "
$ws.Range("D18").Value = "'2.442.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("D19").Value = "'7.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D20").Value = "'10.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.16%  "
$ws.Range("D21").Value = "'328.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("E22").Value = "  +1.27%  "
$ws.Range("D23").Value = "'2.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.72%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'0.0₆0674"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +140.48%  "
$ws.Range("D26").Value = "'65.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("D27").Value = "'622.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.33%  "
$ws.Range("E28").Value = "  +11.46%  "
$ws.Range("D29").Value = "'8.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.25%  "
$ws.Range("D30").Value = "'0.0₃0984"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.51%  "
$ws.Range("D31").Value = "'2.560.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("D32").Value = "'8.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("D33").Value = "'1.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.12%  "
$ws.Range("D34").Value = "'0.140"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.11%  "
$ws.Range("D35").Value = "'1.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("D36").Value = "'1.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +3.73%  "
$ws.Range("D39").Value = "'0.374"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("B40").Value = "'RenderToken"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'5.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.45%  "
$ws.Range("B41").Value = "'Monero"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'151.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").Value = "'18.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.31%  "
$ws.Range("D43").Value = "'2.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.95%  "
$ws.Range("D44").Value = "'1.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.94%  "
$ws.Range("D46").Value = "'0.867"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +126.45%  "
$ws.Range("D47").Value = "'15.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +29.46%  "
$ws.Range("D48").Value = "'144.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("D49").Value = "'3.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").Value = "'20.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.49%  "
$ws.Range("D51").Value = "'0.599"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.78%  "
